# ===========================================================================
# Apply the "adding averages and more checks" update:
#   - Training Dashboard: refresh PERIOD TO EXPIRE (H) and LAST UPDATE (I)
#     figures for every training row (recomputed against a newer "as of"
#     date of 16-Sep-2025 instead of 08-Sep-2025).
#   - Exam Dashboard: narrow the COMMENTS column and replace the old
#     per-row comments with a uniform "date is valid" remark now that the
#     exam-date check has been re-run.
#   - Header rows on both sheets get a white font colour so the bold
#     header text reads correctly against the dark blue fill.
# ===========================================================================

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Training Dashboard")
$ws2 = $wb.Worksheets.Item("Exam Dashboard")

# ---------------------------------------------------------------------
# Helper: write a literal text value into a cell without letting Excel's
# auto-recognition turn a date-looking string (e.g. "16-Sep-2025") into a
# real date serial number / date-formatted cell. We build the text via a
# formula first (formulas are not re-interpreted the way typed values
# are) and then flatten it back down to a plain value with Paste Special
# so the cell ends up as an ordinary text value, keeping its existing
# number format / style untouched.
# ---------------------------------------------------------------------
function Set-LiteralText {
    param($range, [string]$text)

    $escaped = $text.Replace('"', '""')
    $range.Formula = '="' + $escaped + '"'
    $range.Copy() | Out-Null
    $range.PasteSpecial(-4163) | Out-Null   # xlPasteValues
}

# --- Training Dashboard: updated PERIOD TO EXPIRE / LAST UPDATE values ---

$ws1.Range("H3").Value = 265
Set-LiteralText $ws1.Range("I3") "16-Sep-2025"

$ws1.Range("H4").Value = 394
Set-LiteralText $ws1.Range("I4") "16-Sep-2025"

$ws1.Range("H5").Value = 413
Set-LiteralText $ws1.Range("I5") "16-Sep-2025"

$ws1.Range("H6").Value = 413
Set-LiteralText $ws1.Range("I6") "16-Sep-2025"

$ws1.Range("H7").Value = 350
Set-LiteralText $ws1.Range("I7") "16-Sep-2025"

$ws1.Range("H8").Value = -103
Set-LiteralText $ws1.Range("I8") "16-Sep-2025"

$ws1.Range("H9").Value = -328
Set-LiteralText $ws1.Range("I9") "16-Sep-2025"

$ws1.Range("H10").Value = 155
Set-LiteralText $ws1.Range("I10") "16-Sep-2025"

# --- Exam Dashboard: narrower COMMENTS column + refreshed remarks -------

$ws2.Columns.Item(5).ColumnWidth = 14.14   # renders as width 15 in the xlsx

Set-LiteralText $ws2.Range("E3") "date is valid"
Set-LiteralText $ws2.Range("E4") "date is valid"
Set-LiteralText $ws2.Range("E5") "date is valid"

# --- Header rows: bold white text on the existing dark-blue fill --------

$ws1.Range("A2:K2").Font.Color = 16777215
$ws2.Range("A2:G2").Font.Color = 16777215
